$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style) from H1 onto the new I1/J1 header cells,
# then set their text values.
$ws.Range("H1").Copy()
$ws.Range("I1").PasteSpecial(-4122)
$ws.Range("J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new data columns:
#   I (rows 2-8) is always 1
#   J (rows 2-8) mirrors the value already present in column H
for ($r = 2; $r -le 8; $r++) {
    $hVal = $ws.Cells.Item($r, 8).Value2
    $ws.Cells.Item($r, 9).Value = 1
    $ws.Cells.Item($r, 10).Value = $hVal
}
